$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D43").Value = "1969 - 1986"
$ws.Range("D46:D48").Value = "1969 - 1986"
$ws.Range("D53").Value = "1969 - 1986"
$ws.Range("D63").Value = "1969 - 1986"
$ws.Range("D65").Value = "1969 - 1986"
$ws.Range("D67").Value = "1969 - 1986"
$ws.Range("D86:D87").Value = "1969 - 1986"
$ws.Range("D93").Value = "1969 - 1986"
$ws.Range("D96").Value = "1969 - 1986"

$ws.Range("D15").Value = "1987 - 2004"
$ws.Range("D21").Value = "1987 - 2004"
$ws.Range("D52").Value = "1987 - 2004"
$ws.Range("D68").Value = "1987 - 2004"
$ws.Range("D71").Value = "1987 - 2004"
$ws.Range("D89").Value = "1987 - 2004"

$ws.Range("D2:D9").Value = "2005 - 2022"
$ws.Range("D11:D14").Value = "2005 - 2022"
$ws.Range("D16:D20").Value = "2005 - 2022"
$ws.Range("D22:D31").Value = "2005 - 2022"
$ws.Range("D33:D42").Value = "2005 - 2022"
$ws.Range("D44").Value = "2005 - 2022"
$ws.Range("D49:D51").Value = "2005 - 2022"
$ws.Range("D54:D62").Value = "2005 - 2022"
$ws.Range("D64").Value = "2005 - 2022"
$ws.Range("D69:D70").Value = "2005 - 2022"
$ws.Range("D72:D84").Value = "2005 - 2022"
$ws.Range("D88").Value = "2005 - 2022"
$ws.Range("D90:D92").Value = "2005 - 2022"

